$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 501, pushing existing rows 501-620 down to 502-621.
$ws.Rows(501).Insert()

# Populate the newly inserted row 501 with the new record.
$ws.Cells.Item(501, 1).Value = 7
$ws.Cells.Item(501, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(501, 3).Value = "Ñuble"
$ws.Cells.Item(501, 4).Value = 45204
$ws.Cells.Item(501, 5).Value = 16
$ws.Cells.Item(501, 6).Value = 100112023
$ws.Cells.Item(501, 7).Value = "Brócoli"
$ws.Cells.Item(501, 8).Value = "Sin especificar"
$ws.Cells.Item(501, 9).Value = "Primera"
$ws.Cells.Item(501, 10).Value = 200
$ws.Cells.Item(501, 11).Value = 1200
$ws.Cells.Item(501, 12).Value = 1200
$ws.Cells.Item(501, 13).Value = 1200
$ws.Cells.Item(501, 14).Value = "$/unidad"
$ws.Cells.Item(501, 15).Value = "Región del Maule"
$ws.Cells.Item(501, 16).Value = 1200
$ws.Cells.Item(501, 17).Value = 1
$ws.Cells.Item(501, 18).Value = "Hortaliza"
